$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "324.34") must be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# values and the distinctive formatting ("1.090" -> 1.09, "7.800" -> 7.8, etc.)
# gets lost. We flip the format to Text, write the values, then restore the
# default "Normal" style so no extra formatting is left behind on the cells.
$numericLookingCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D17", "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "324.34"
$ws.Range("D6").Value = "1.001"
$ws.Range("D7").Value = "0.4562"
$ws.Range("D8").Value = "0.3553"
$ws.Range("D9").Value = "0.07470"
$ws.Range("D10").Value = "41.54"
$ws.Range("D11").Value = "1.088"
$ws.Range("D13").Value = "20.79"
$ws.Range("D14").Value = "6.019"
$ws.Range("D17").Value = "94.03"
$ws.Range("D19").Value = "0.06402"
$ws.Range("D21").Value = "17.13"
$ws.Range("D22").Value = "5.747"
$ws.Range("D24").Value = "11.21"
$ws.Range("D25").Value = "2.083"
$ws.Range("D26").Value = "165.67"
$ws.Range("D27").Value = "20.14"
$ws.Range("D29").Value = "2.135"
$ws.Range("D30").Value = "125.76"
$ws.Range("D31").Value = "1.090"
$ws.Range("D33").Value = "3.661"
$ws.Range("D34").Value = "5.535"
$ws.Range("D35").Value = "11.73"
$ws.Range("D37").Value = "0.2094"
$ws.Range("D38").Value = "0.06018"
$ws.Range("D39").Value = "0.6303"
$ws.Range("D40").Value = "4.922"
$ws.Range("D41").Value = "1.180"
$ws.Range("D42").Value = "1.387"
$ws.Range("D43").Value = "7.800"
$ws.Range("D44").Value = "13.22"
$ws.Range("D46").Value = "0.5862"
$ws.Range("D47").Value = "122.10"
$ws.Range("D48").Value = "1.933"
$ws.Range("D49").Value = "0.06893"
$ws.Range("D51").Value = "72.11"

foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining cells are safe to set directly (their text does not parse as a number).
$ws.Range("D2").Value = "27.548.61"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.754.59"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("E8").Value = "  -1.94%  "
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D16").Value = "1.759.93"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("D23").Value = "27.603.29"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  +1.94%  "
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "1.958.06"
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("E51").Value = "  -0.47%  "
